$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (shared strings): month names shift forward by one month.
# Assigning strings like "March 2024" directly to .Value gets auto-converted to a
# date serial by Excel's smart-typing. To avoid that (and avoid leaving behind a
# custom number format on the cell), build the text via a temporary formula cell,
# copy it, and paste-special the resulting value back into the target cell.
$ws.Range("I1").Formula = '="March 2024"'
$ws.Range("I1").Copy()
$ws.Range("A1").PasteSpecial(-4163)

$ws.Range("I1").Formula = '="April 2024"'
$ws.Range("I1").Copy()
$ws.Range("G1").PasteSpecial(-4163)

$ws.Range("I1").Value = ""

# Update row 2 numeric values
$ws.Range("A2").Value = 1.658
$ws.Range("B2").Value = 0.233
$ws.Range("C2").Value = 0.061
$ws.Range("D2").Value = -0.038
$ws.Range("E2").Value = -0.007
$ws.Range("F2").Value = -0.304
$ws.Range("G2").Value = 1.602
